$d = $word.ActiveDocument

# Move the _GoBack bookmark from its old location (end of the "one hot
# encoding." paragraph) to the last paragraph of the document body (just
# before the final sectPr). Bookmark names are unique, so adding a new
# "_GoBack" bookmark removes the old one automatically.
$paras = $d.Paragraphs
$lastPara = $paras.Item($paras.Count)
$d.Bookmarks.Add("_GoBack", $lastPara.Range)

# Mark "Normal (Web)" and "_Style 8" as Quick Styles (adds <w:qFormat/>).
$d.Styles("Normal (Web)").QuickStyle = $true
$d.Styles("_Style 8").QuickStyle = $true
